# Auto-upload VRF Excel file
#
# Adds a new worksheet named "test" at the very end of the workbook.
# The sheet uses the same "Outdoor/Indoor Model/Quantity/Serial(s)"
# header layout as every other sheet in this workbook, so the new tab
# is created by duplicating the last existing sheet (which already
# carries the right header text/format/page setup) and then replacing
# its body with a couple of rows of test data.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate the last sheet (keeps header formatting / page setup) and
# place the copy right after it, i.e. at the very end of the tab strip.
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "test"

# Wipe the copied sheet's old data rows, keeping only the header row.
$ws.Cells.ClearContents()

# Re-write the header row (Outdoor/Indoor Model/Quantity/Serial(s)).
$ws.Range("A1").Value = "Outdoor Model"
$ws.Range("B1").Value = "Outdoor Quantity"
$ws.Range("C1").Value = "Outdoor Serial(s)"
$ws.Range("D1").Value = "Indoor Model"
$ws.Range("E1").Value = "Indoor Quantity"
$ws.Range("F1").Value = "Indoor Serial(s)"

# Row 2 test data.
$ws.Range("A2").Value = "sf"
$ws.Range("B2").Value = 0
# Remaining row-2 cells were touched (typed into) but left blank.
$ws.Range("C2").Font.Bold = $false
$ws.Range("D2").Font.Bold = $false
$ws.Range("E2").Font.Bold = $false
$ws.Range("F2").Font.Bold = $false

# Row 3 test data.
$ws.Range("A3").Font.Bold = $false
$ws.Range("B3").Font.Bold = $false
$ws.Range("C3").Font.Bold = $false
$ws.Range("D3").Value = "dtgdse"
$ws.Range("E3").Value = 0
$ws.Range("F3").Font.Bold = $false

$ws.Range("A1").Select()
